$wb = $excel.ActiveWorkbook

# --- Update the conversion text on "Hoja1" ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsHoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 3.12 = 11621.74 pesos`n✅ 11621.74 pesos = 3.09 = 953.54 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Update the rate figures on "tasas" ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 320.52
$wsTasas.Range("O10").Value = 3725
$wsTasas.Range("N12").Value = 3760.01
$wsTasas.Range("O12").Value = 308.5
